# Klipper US.xlsx edit: add "US 8" user-story sheet, rename "US7" -> "US7 "
$wb = $excel.ActiveWorkbook

# 1. Rename the "US7" tab to "US7 " (trailing space, as in target workbook).
$ws7 = $wb.Worksheets.Item("US7")
$ws7.Name = "US7 "

# 2. Create the new "US8" sheet by copying the "US7 " sheet (same layout/
#    styling as the other "US N" task sheets) and placing it right before
#    the "Issues" sheet.
$issues = $wb.Worksheets.Item("Issues")
$ws7.Copy($issues)
$ws8 = $wb.Worksheets.Item("US7  (2)")
$ws8.Name = "US8"

# 3. Trim the copied sheet down to just the rows this user story needs
#    (header @ row3, blank @ row4, column headers @ row5, 3 task rows).
$ws8.Rows("9:19").Delete()

# 4. Fill in the US 8 content.
$ws8.Range("B3").Value = "US 8"
$ws8.Range("C3").Value = "Show Total hrs for the selected data"

$ws8.Range("B6").Value = 1
$ws8.Range("C6").Value = "Calculate Total working hours and total deficit/overtime hours for specified date range"
$ws8.Range("D6").Value = 2
$ws8.Range("E6").Value = "Sanket"
$ws8.Range("F6").Value = "To do"

$ws8.Range("B7").Value = 2
$ws8.Range("C7").Value = "UI - show total working hours and deficit/overtime hours "
$ws8.Range("D7").Value = 2
$ws8.Range("E7").Value = "Shweta"
$ws8.Range("F7").Value = "To do"

$ws8.Range("B8").Value = 3
$ws8.Range("C8").Value = "write test cases"
$ws8.Range("D8").Value = 2
$ws8.Range("E8").Value = "Sidhdesh"
$ws8.Range("F8").Value = "To do"

# 5. Widen column C so the longer task descriptions are readable, and make
#    the new sheet the active tab/selection (matches the authored workbook).
$ws8.Columns("C:C").AutoFit()
$ws8.Activate()
[void]$ws8.Range("C7").Select()
